$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.017.87"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "2.600.34"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.27"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.00"
$ws.Range("E6").Value = "  -2.47%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("D9").Value = "2.600.71"
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.128"
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.342"
$ws.Range("E13").Value = "  -3.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.08"
$ws.Range("E14").Value = "  -2.92%  "
$ws.Range("D15").Value = "3.074.45"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("E16").Value = "  -3.17%  "
$ws.Range("D17").Value = "66.993.42"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").Value = "2.605.88"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "366.89"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.01"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.36"
$ws.Range("E21").Value = "  -3.44%  "
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.73"
$ws.Range("E23").Value = "  -4.81%  "
$ws.Range("E24").Value = "  -3.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.41"
$ws.Range("E25").Value = "  +4.87%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.90"
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").Value = "2.735.60"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "581.93"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").Value = "0.0₃0982"
$ws.Range("E31").Value = "  -7.07%  "
$ws.Range("E32").Value = "  -5.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.64"
$ws.Range("E33").Value = "  -3.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  -4.28%  "
$ws.Range("E37").Value = "  -2.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.60"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.99"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.84"
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.21"
$ws.Range("E42").Value = "  -3.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.55"
$ws.Range("E43").Value = "  -4.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.09"
$ws.Range("E44").Value = "  +4.14%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.53"
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("D47").Value = "0.0₆0285"
$ws.Range("E47").Value = "  -2.17%  "
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0779"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.67"
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.29"
$ws.Range("E51").Value = "  +1.51%  "
